# "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# This "Estado de Cuenta" sheet has a small table of worker contribution
# periods (rows 16-24, period codes 2412, 2501-2508) followed by a
# signature block. The edit:
#   1. Adds a new period row (2509) for the same worker, right after the
#      current last row (2508) — pushing the blank spacer rows and the
#      signature block down by one row.
#   2. Updates the "VALOR MORA" total and the "Cant. Periodos" count to
#      reflect the newly added period.
#   3. Centers the period-code column (E) in the table, matching the
#      refreshed template formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new table row right after the current last data row (24) ---
$ws.Rows("25:25").Insert()

# The freshly inserted row starts out with a generic/blank style. Give it
# the closing "bottom of table" look that the old row 24 had...
$ws.Range("B24:J24").Copy()
$ws.Range("B25:J25").PasteSpecial(-4122)

# ...and demote the old row 24 (2508) to a regular interior row, matching
# the style used by rows 16-23.
$ws.Range("B23:J23").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Populate the new row with the new period for the same worker ---
$ws.Range("B25").Value2 = "CC"
$ws.Range("C25").Value2 = "73578753"
$ws.Range("D25").Value2 = "DARWIN ALEXANDER MERCADO VEGA"
$ws.Range("E25").Value2 = "2509"
$ws.Range("F25").Value2 = 52000
$ws.Range("G25").Value2 = 1300000
$ws.Range("H25").Value2 = ""
$ws.Range("I25").Value2 = ""
$ws.Range("J25").Value2 = ""

# --- 3. Refresh the summary figures ---
$ws.Range("E11").Value2 = 481867
$ws.Range("F13").Value2 = 10

# --- 4. Center the "Periodo Mora" column across the whole table ---
$ws.Range("E16:E25").HorizontalAlignment = -4108

Write-Host "Edit complete"
